# Auto-generated Excel COM script to apply the odds update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.73
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 2.38
$ws.Range("X2").Value = 23
$ws.Range("Z2").Value = 51
$ws.Range("AD2").Value = 7.5
$ws.Range("AI2").Value = 8
$ws.Range("AK2").Value = 13
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 26
$ws.Range("AS2").Value = 251
$ws.Range("AX2").Value = 9
$ws.Range("G3").Value = 1.3
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 1.8
$ws.Range("K3").Value = 2.5
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 4.33
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.4
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 6.5
$ws.Range("Y3").Value = 9
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 13
$ws.Range("AD3").Value = 9.5
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 23
$ws.Range("AM3").Value = 67
$ws.Range("AN3").Value = 3.25
$ws.Range("AO3").Value = 6
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 15
$ws.Range("AR3").Value = 41
$ws.Range("AS3").Value = 126
$ws.Range("AT3").Value = 3.4
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 9.5
$ws.Range("BA3").Value = 201
$ws.Range("BB3").Value = 351
$ws.Range("G4").Value = 2.1
$ws.Range("I4").Value = 3.4
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.85
$ws.Range("X4").Value = 10
$ws.Range("Z4").Value = 19
$ws.Range("AB4").Value = 26
$ws.Range("AJ4").Value = 12
$ws.Range("AO4").Value = 12
$ws.Range("AX4").Value = 19
$ws.Range("G5").Value = 2.35
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.2
$ws.Range("L5").Value = 4
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.75
$ws.Range("X5").Value = 10
$ws.Range("AI5").Value = 15
$ws.Range("AW5").Value = 5
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 3.2
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.95
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("AB6").Value = 23
$ws.Range("AG6").Value = 151
$ws.Range("AH6").Value = 11
$ws.Range("AT6").Value = 3
$ws.Range("BC6").Value = 501
$ws.Range("G7").Value = 2.9
$ws.Range("I7").Value = 2.45
$ws.Range("L7").Value = 3.4
$ws.Range("X7").Value = 13
$ws.Range("Z7").Value = 34
$ws.Range("AA7").Value = 29
$ws.Range("AJ7").Value = 10
$ws.Range("AK7").Value = 23
$ws.Range("AO7").Value = 19
$ws.Range("AP7").Value = 34
$ws.Range("AQ7").Value = 67
$ws.Range("AW7").Value = 4.33
$ws.Range("AX7").Value = 15
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.67
$ws.Range("H9").Value = 3.25
$ws.Range("N9").Value = 8.5
$ws.Range("AL9").Value = 23
$ws.Range("AY9").Value = 26
$ws.Range("A11").Value = '2FstQtPr'
$ws.Range("E11").Value = 'PSIS Semarang'
$ws.Range("F11").Value = 'Bali United'
$ws.Range("G11").Value = 3.3
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 2.15
$ws.Range("J11").Value = 3.8
$ws.Range("K11").Value = 2.05
$ws.Range("L11").Value = 2.7
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 8.1
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 2.77
$ws.Range("Q11").Value = 1.98
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 1.39
$ws.Range("T11").Value = 2.55
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.85
$ws.Range("W11").Value = 9.25
$ws.Range("X11").Value = 17
$ws.Range("Y11").Value = 11.5
$ws.Range("Z11").Value = 45
$ws.Range("AA11").Value = 32
$ws.Range("AB11").Value = 40
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 6
$ws.Range("AE11").Value = 14.5
$ws.Range("AF11").Value = 70
$ws.Range("AG11").Value = 600
$ws.Range("AH11").Value = 7.1
$ws.Range("AI11").Value = 10.25
$ws.Range("AJ11").Value = 8.75
$ws.Range("AK11").Value = 21
$ws.Range("AL11").Value = 18
$ws.Range("AM11").Value = 29
$ws.Range("AN11").Value = 5.2
$ws.Range("AO11").Value = 18.5
$ws.Range("AP11").Value = 24
$ws.Range("AQ11").Value = 90
$ws.Range("AR11").Value = 120
$ws.Range("AS11").Value = 300
$ws.Range("AT11").Value = 2.52
$ws.Range("AU11").Value = 6.8
$ws.Range("AV11").Value = 60
$ws.Range("AW11").Value = 4
$ws.Range("AX11").Value = 11
$ws.Range("AY11").Value = 18.5
$ws.Range("AZ11").Value = 40
$ws.Range("BA11").Value = 70
$ws.Range("BB11").Value = 250
$ws.Range("A12").Value = 'bqG8rkuA'
$ws.Range("E12").Value = 'Persebaya'
$ws.Range("F12").Value = 'Persik Kediri'
$ws.Range("G12").Value = 1.5
$ws.Range("H12").Value = 3.95
$ws.Range("I12").Value = 5.8
$ws.Range("J12").Value = 2.05
$ws.Range("K12").Value = 2.18
$ws.Range("L12").Value = 5.7
$ws.Range("M12").Value = 1.02
$ws.Range("N12").Value = 7.5
$ws.Range("O12").Value = 1.27
$ws.Range("P12").Value = 3.1
$ws.Range("Q12").Value = 1.8
$ws.Range("R12").Value = 1.82
$ws.Range("S12").Value = 1.38
$ws.Range("T12").Value = 2.6
$ws.Range("U12").Value = 1.9
$ws.Range("V12").Value = 1.72
$ws.Range("W12").Value = 6.4
$ws.Range("X12").Value = 6.8
$ws.Range("Y12").Value = 8.25
$ws.Range("Z12").Value = 10.25
$ws.Range("AA12").Value = 12.5
$ws.Range("AB12").Value = 29
$ws.Range("AC12").Value = 10.25
$ws.Range("AD12").Value = 7.8
$ws.Range("AE12").Value = 19
$ws.Range("AF12").Value = 100
$ws.Range("AG12").Value = 800
$ws.Range("AH12").Value = 14.5
$ws.Range("AI12").Value = 35
$ws.Range("AJ12").Value = 19
$ws.Range("AK12").Value = 120
$ws.Range("AL12").Value = 65
$ws.Range("AM12").Value = 65
$ws.Range("AN12").Value = 3.25
$ws.Range("AO12").Value = 7.1
$ws.Range("AP12").Value = 18
$ws.Range("AQ12").Value = 22
$ws.Range("AR12").Value = 55
$ws.Range("AS12").Value = 250
$ws.Range("AT12").Value = 2.55
$ws.Range("AU12").Value = 8.25
$ws.Range("AV12").Value = 90
$ws.Range("AW12").Value = 7.1
$ws.Range("AX12").Value = 35
$ws.Range("AY12").Value = 40
$ws.Range("AZ12").Value = 250
$ws.Range("BA12").Value = 300
$ws.Range("BB12").Value = 500
